$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @{
    2  = @(3, 1, 67.77251700000001, 203.317551, 0.4079637943863715, 0.4079637943863715, 3, 1, 2.672264666666667, 8.016794000000001, 0.06772620019093417, 0.06772620019093417, 181.106102550166, 1629.954922951494, 0.0276298376092645, 0.0276298376092645)
    3  = @(3, 1, 67.77251700000001, 203.317551, 0.4079637943863715, 0.4079637943863715, 3, 1, 27.41886, 82.25658, 0.6949069171668364, 0.6949069171668364, 1858.24515547062, 16724.20639923558, 0.2834968626727185, 0.2834968626727185)
    4  = @(3, 1, 67.77251700000001, 203.317551, 0.4079637943863715, 0.4079637943863715, 3, 1, 9.365757, 28.097271, 0.2373668826422294, 0.2373668826422294, 634.7409255003691, 5712.668329503321, 0.09683709410438845, 0.09683709410438845)
    5  = @(3, 1, 60.97760633333333, 182.932819, 0.3670611149405164, 0.3670611149405164, 3, 1, 2.672264666666667, 8.016794000000001, 0.06772620019093417, 0.06772620019093417, 162.9483028624763, 1466.534725762286, 0.02485965455276891, 0.02485965455276891)
    6  = @(3, 1, 60.97760633333333, 182.932819, 0.3670611149405164, 0.3670611149405164, 3, 1, 27.41886, 82.25658, 0.6949069171668364, 0.6949069171668364, 1671.93645118878, 15047.42806069902, 0.255073307795136, 0.255073307795136)
    7  = @(3, 1, 60.97760633333333, 182.932819, 0.3670611149405164, 0.3670611149405164, 3, 1, 9.365757, 28.097271, 0.2373668826422294, 0.2373668826422294, 571.1014433596611, 5139.912990236949, 0.08712815259261143, 0.08712815259261142)
    8  = @(3, 1, 37.37372866666667, 112.121186, 0.2249750906731122, 0.2249750906731122, 3, 1, 2.672264666666667, 8.016794000000001, 0.06772620019093417, 0.06772620019093417, 99.87249457752047, 898.8524511976842, 0.01523670802890076, 0.01523670802890076)
    9  = @(3, 1, 37.37372866666667, 112.121186, 0.2249750906731122, 0.2249750906731122, 3, 1, 27.41886, 82.25658, 0.6949069171668364, 0.6949069171668364, 1024.74503398932, 9222.705305903881, 0.1563367466989819, 0.1563367466989819)
    10 = @(3, 1, 37.37372866666667, 112.121186, 0.2249750906731122, 0.2249750906731122, 3, 1, 9.365757, 28.097271, 0.2373668826422294, 0.2373668826422294, 350.033260875934, 3150.299347883406, 0.05340163594522954, 0.05340163594522954)
}

# Columns E(5) through T(20)
$startCol = 5

foreach ($row in $data.Keys) {
    $values = $data[$row]
    for ($i = 0; $i -lt $values.Length; $i++) {
        $col = $startCol + $i
        $ws.Cells.Item($row, $col).Value = $values[$i]
    }
}
